# Product en sprint backlog consistent geschreven, conform de Scrumgids
# (behalve hoofdletters). Hernoem de "BL" (backlog) afkorting naar "PB"
# (product backlog) op de dia met de legenda/afkortingen.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    if (-not $sh.HasTextFrame) {
        continue
    }
    if (-not $sh.TextFrame.HasText) {
        continue
    }

    $tr = $sh.TextFrame.TextRange
    $full = $tr.Text

    # Case 1: the folded-corner ("Ezelsoor") legend box whose whole text is
    # exactly the abbreviation "BL".
    if ($full -eq "BL") {
        $tr.Text = "PB"
        continue
    }

    # Case 2: the legend textbox ("Tekstvak 109") that lists every
    # abbreviation; one paragraph/run starts with "BL<TAB>product ".
    $needle = "BL" + [char]9 + "product "
    $idx = $full.IndexOf($needle)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $needle.Length)
        $sub.Text = "PB" + [char]9 + "product "
    }
}

# Best-effort: the canonical diff also bumps the cached text of the two
# "datetimeFigureOut" auto-date fields (handout master & notes master) from
# 03-02-2022 to 28-02-2025 -- this is simply PowerPoint re-stamping the
# "last saved" automatic date, not a deliberate content edit, so it is
# applied here defensively in case the host allows rewriting it.
$targetDate = "28-02-2025"

$hm = $p.HandoutMaster
for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
    $sh = $hm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "03-02-2022") {
            $sh.TextFrame.TextRange.Text = $targetDate
        }
    }
}

$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "03-02-2022") {
            $sh.TextFrame.TextRange.Text = $targetDate
        }
    }
}
